$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.620.64'
$ws.Range('E2').Value = '  +1.80%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.994.96'
$ws.Range('E3').Value = '  +3.20%  '
$ws.Range('E4').Value = '  +0.32%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '381.74'
$ws.Range('E5').Value = '  +5.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '106.19'
$ws.Range('E6').Value = '  +1.86%  '
$ws.Range('E7').Value = '  +1.31%  '
$ws.Range('E8').Value = '  -2.66%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.600'
$ws.Range('E9').Value = '  +1.93%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.54'
$ws.Range('E10').Value = '  +1.75%  '
$ws.Range('E11').Value = '  +0.89%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0846'
$ws.Range('E12').Value = '  +1.27%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.68'
$ws.Range('E13').Value = '  +1.74%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.463.23'
$ws.Range('E14').Value = '  +3.35%  '
$ws.Range('E15').Value = '  +2.42%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.992.73'
$ws.Range('E16').Value = '  +3.25%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.982'
$ws.Range('E17').Value = '  +3.43%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '51.674.13'
$ws.Range('E18').Value = '  +1.88%  '
$ws.Range('E19').Value = '  +3.26%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.45'
$ws.Range('E20').Value = '  +2.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.10'
$ws.Range('E21').Value = '  +0.93%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0965'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.24'
$ws.Range('E23').Value = '  +1.76%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '263.43'
$ws.Range('E24').Value = '  +1.34%  '
$ws.Range('E25').Value = '  +5.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.171'
$ws.Range('E26').Value = '  -1.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.22'
$ws.Range('E27').Value = '  +18.36%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.48'
$ws.Range('E28').Value = '  +1.84%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '26.12'
$ws.Range('E29').Value = '  +1.03%  '
$ws.Range('B30').Value = 'Dai'
$ws.Range('C30').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('E31').Value = '  +6.44%  '
$ws.Range('E32').Value = '  -0.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '35.25'
$ws.Range('E33').Value = '  +1.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.10'
$ws.Range('E34').Value = '  -1.90%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '51.33'
$ws.Range('E35').Value = '  +2.13%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0460'
$ws.Range('E36').Value = '  +9.43%  '
$ws.Range('E37').Value = '  +0.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.12'
$ws.Range('E38').Value = '  -0.04%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '17.57'
$ws.Range('E39').Value = '  +3.64%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.63'
$ws.Range('E40').Value = '  -5.15%  '
$ws.Range('E41').Value = '  +0.50%  '
$ws.Range('E42').Value = '  +2.93%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '124.19'
$ws.Range('E43').Value = '  +6.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '22.38'
$ws.Range('E44').Value = '  -16.73%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.282'
$ws.Range('E45').Value = '  +20.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.08'
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('E47').Value = '  +5.24%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.053.16'
$ws.Range('E48').Value = '  -0.07%  '
$ws.Range('E49').Value = '  +2.38%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0360'
$ws.Range('E50').Value = '  +11.94%  '
$ws.Range('E51').Value = '  +4.39%  '
